# ldlc_suivi_smartphones.xlsx - "Update LDLC prices history"
#
# A new timestamped price-check column is inserted right before the
# existing "nom" / "url_produit" columns (previously CR / CS, now pushed
# out to CS / CT). The new column is populated with the latest price
# snapshot (copied from the preceding price column, CQ before the
# insert) for every product row that already has a numeric price there;
# rows without a price yet are left blank, matching the existing pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column CR = 96 (1-based). Inserting here shifts the old CR ("nom") to
# CS and the old CS ("url_produit") to CT, carrying styles/values along.
$newCol = 96      # CR
$sourceCol = 95   # CQ (the last existing price-snapshot column)

$ws.Columns("CR").Insert()

# Header row: new timestamp label for the freshly inserted column.
$ws.Cells.Item(1, $newCol).Value2 = "2026-02-01 01:09:08"

# Data rows: rows 2-80 already carry a numeric price in CQ; mirror it
# into the new CR column. Rows 81-206 have no price yet in CQ, so the
# new CR cell is left empty, same as the source.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, $sourceCol).Value2
    if ($v -ne $null -and $v -ne "") {
        $ws.Cells.Item($r, $newCol).Value2 = $v
    }
}
